$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5976617336273193
$ws.Range("E2").Value = 1806.468575708612
$ws.Range("F2").Value = 0.08560086259044505
$ws.Range("G2").Value = 0.06574912211457055
$ws.Range("H2").Value = 0.05733266186615138
$ws.Range("I2").Value = 0.05066433169342888
$ws.Range("J2").Value = 0.0480582989953702
$ws.Range("K2").Value = 0.04518442706794643
$ws.Range("L2").Value = 0.04383405651500173
$ws.Range("M2").Value = 0.04139457088581598
$ws.Range("N2").Value = 0.0407818970491105
$ws.Range("O2").Value = 0.03998739670879708
$ws.Range("P2").Value = 0.03936233586700081
$ws.Range("Q2").Value = 0.037936642660431
$ws.Range("R2").Value = 0.03700533369581052
$ws.Range("S2").Value = 0.03661382972024962
$ws.Range("T2").Value = 0.03617042805258851
$ws.Range("U2").Value = 0.03617042805258851
$ws.Range("V2").Value = 0.03575799179583767
$ws.Range("W2").Value = 0.03560555439549738
$ws.Range("X2").Value = 0.03541846365061712
$ws.Range("Y2").Value = 0.03521381239198073
$ws.Range("C3").Value = 0.6562473773956299
$ws.Range("E3").Value = 1807.665806068271
$ws.Range("F3").Value = 0.08926126878392274
$ws.Range("G3").Value = 0.06700466454246336
$ws.Range("H3").Value = 0.05952814514996359
$ws.Range("I3").Value = 0.04894187871438552
$ws.Range("J3").Value = 0.04894187871438552
$ws.Range("K3").Value = 0.04676211193500163
$ws.Range("L3").Value = 0.0429695971775418
$ws.Range("M3").Value = 0.04216994795310181
$ws.Range("N3").Value = 0.04083413420651079
$ws.Range("O3").Value = 0.03932713636711409
$ws.Range("P3").Value = 0.03871332503778702
$ws.Range("Q3").Value = 0.03854840727114529
$ws.Range("R3").Value = 0.03785455783398946
$ws.Range("S3").Value = 0.03715921629398585
$ws.Range("T3").Value = 0.03711638117820813
$ws.Range("U3").Value = 0.03585138749227134
$ws.Range("V3").Value = 0.03585138749227134
$ws.Range("W3").Value = 0.03561642812237189
$ws.Range("X3").Value = 0.03561642812237189
$ws.Range("Y3").Value = 0.03523715021575577
$ws.Range("C4").Value = 0.6241059303283691
$ws.Range("E4").Value = 1773.733536374357
$ws.Range("F4").Value = 0.08537155427990441
$ws.Range("G4").Value = 0.06824074674906187
$ws.Range("H4").Value = 0.06101614032988132
$ws.Range("I4").Value = 0.05439879326683741
$ws.Range("J4").Value = 0.04781315108862522
$ws.Range("K4").Value = 0.04651609695628686
$ws.Range("L4").Value = 0.04381048617728529
$ws.Range("M4").Value = 0.04223914149851184
$ws.Range("N4").Value = 0.03929323660679783
$ws.Range("O4").Value = 0.03892568734778481
$ws.Range("P4").Value = 0.03875254816293498
$ws.Range("Q4").Value = 0.03759809293237855
$ws.Range("R4").Value = 0.03717889621805095
$ws.Range("S4").Value = 0.03674632910281461
$ws.Range("T4").Value = 0.0365726846530563
$ws.Range("U4").Value = 0.03613867491293095
$ws.Range("V4").Value = 0.03549840654224003
$ws.Range("W4").Value = 0.03522482179773609
$ws.Range("X4").Value = 0.03486629640336285
$ws.Range("Y4").Value = 0.03457570246343775
$ws.Range("C5").Value = 0.5312480926513672
$ws.Range("E5").Value = 1761.151153211811
$ws.Range("F5").Value = 0.08482318001874259
$ws.Range("G5").Value = 0.06336794811751684
$ws.Range("H5").Value = 0.05808245563486555
$ws.Range("I5").Value = 0.05255624623779582
$ws.Range("J5").Value = 0.04603112309988559
$ws.Range("K5").Value = 0.04343216200270867
$ws.Range("L5").Value = 0.04293393570594766
$ws.Range("M5").Value = 0.04157958151993055
$ws.Range("N5").Value = 0.04108009079121559
$ws.Range("O5").Value = 0.040214537238629
$ws.Range("P5").Value = 0.03948474539108084
$ws.Range("Q5").Value = 0.03783978995028617
$ws.Range("R5").Value = 0.03745844160536983
$ws.Range("S5").Value = 0.03643150565684686
$ws.Range("T5").Value = 0.03643150565684686
$ws.Range("U5").Value = 0.03546920090280604
$ws.Range("V5").Value = 0.03511271667064466
$ws.Range("W5").Value = 0.0350168999867165
$ws.Range("X5").Value = 0.03463301586213853
$ws.Range("Y5").Value = 0.03433043183648754
$ws.Range("C6").Value = 0.5312614440917969
$ws.Range("E6").Value = 1852.757724778887
$ws.Range("F6").Value = 0.0863679041360383
$ws.Range("G6").Value = 0.06583296023917092
$ws.Range("H6").Value = 0.06000978761518342
$ws.Range("I6").Value = 0.05147168892603441
$ws.Range("J6").Value = 0.04721144784750572
$ws.Range("K6").Value = 0.04391555223857158
$ws.Range("L6").Value = 0.04195187035408634
$ws.Range("M6").Value = 0.04080487438478174
$ws.Range("N6").Value = 0.04080487438478174
$ws.Range("O6").Value = 0.03910406730974397
$ws.Range("P6").Value = 0.03910406730974397
$ws.Range("Q6").Value = 0.03850587367765488
$ws.Range("R6").Value = 0.03829645229366953
$ws.Range("S6").Value = 0.03768462273473294
$ws.Range("T6").Value = 0.03687012163969161
$ws.Range("U6").Value = 0.03687012163969161
$ws.Range("V6").Value = 0.03687012163969161
$ws.Range("W6").Value = 0.0366108320933393
$ws.Range("X6").Value = 0.0366108320933393
$ws.Range("Y6").Value = 0.0361161349859432
$ws.Range("C7").Value = 0.5312504768371582
$ws.Range("E7").Value = 1724.182085578615
$ws.Range("F7").Value = 0.08512857138646447
$ws.Range("G7").Value = 0.06450526370503103
$ws.Range("H7").Value = 0.05521464252854071
$ws.Range("I7").Value = 0.05035078510835539
$ws.Range("J7").Value = 0.04819242687295605
$ws.Range("K7").Value = 0.04502861270157941
$ws.Range("L7").Value = 0.04345230591238183
$ws.Range("M7").Value = 0.04030346733668112
$ws.Range("N7").Value = 0.04030346733668112
$ws.Range("O7").Value = 0.03886844001359686
$ws.Range("P7").Value = 0.03695604488102404
$ws.Range("Q7").Value = 0.03682704807558641
$ws.Range("R7").Value = 0.03654882779917468
$ws.Range("S7").Value = 0.03551882953014283
$ws.Range("T7").Value = 0.03551882953014283
$ws.Range("U7").Value = 0.03458832289149294
$ws.Range("V7").Value = 0.03444975431453014
$ws.Range("W7").Value = 0.03417814335690385
$ws.Range("X7").Value = 0.03407139130590393
$ws.Range("Y7").Value = 0.03360978724324785
$ws.Range("C8").Value = 0.5468764305114746
$ws.Range("E8").Value = 1755.495879773302
$ws.Range("F8").Value = 0.08718237731886827
$ws.Range("G8").Value = 0.06558120120366363
$ws.Range("H8").Value = 0.05871634743644858
$ws.Range("I8").Value = 0.04793055603178414
$ws.Range("J8").Value = 0.04667405391484881
$ws.Range("K8").Value = 0.04537609664418107
$ws.Range("L8").Value = 0.04211558706205529
$ws.Range("M8").Value = 0.04021780178452201
$ws.Range("N8").Value = 0.04017418625259082
$ws.Range("O8").Value = 0.03705305364947834
$ws.Range("P8").Value = 0.03705305364947834
$ws.Range("Q8").Value = 0.03705305364947834
$ws.Range("R8").Value = 0.03655147524809303
$ws.Range("S8").Value = 0.03585718174033449
$ws.Range("T8").Value = 0.03502515297886575
$ws.Range("U8").Value = 0.03476700000337292
$ws.Range("V8").Value = 0.03458727924199133
$ws.Range("W8").Value = 0.03428485331567468
$ws.Range("X8").Value = 0.03428485331567468
$ws.Range("Y8").Value = 0.03422019258817351
$ws.Range("C9").Value = 0.5468649864196777
$ws.Range("E9").Value = 1715.432660156823
$ws.Range("F9").Value = 0.08792296719804281
$ws.Range("G9").Value = 0.06582458834481644
$ws.Range("H9").Value = 0.05893347864143205
$ws.Range("I9").Value = 0.05076377973140533
$ws.Range("J9").Value = 0.04623350212370703
$ws.Range("K9").Value = 0.04411371696620877
$ws.Range("L9").Value = 0.04280404679593963
$ws.Range("M9").Value = 0.03959158036598972
$ws.Range("N9").Value = 0.03823349180385541
$ws.Range("O9").Value = 0.03784569161632855
$ws.Range("P9").Value = 0.03784202451084025
$ws.Range("Q9").Value = 0.03667197913648539
$ws.Range("R9").Value = 0.03585499542409722
$ws.Range("S9").Value = 0.0354022151548756
$ws.Range("T9").Value = 0.03459795445528014
$ws.Range("U9").Value = 0.034336525319059
$ws.Range("V9").Value = 0.0341667191189395
$ws.Range("W9").Value = 0.03370228935135031
$ws.Range("X9").Value = 0.03360355012582299
$ws.Range("Y9").Value = 0.03343923314145852
$ws.Range("C10").Value = 0.5468835830688477
$ws.Range("E10").Value = 1807.825521385472
$ws.Range("F10").Value = 0.08601449778816503
$ws.Range("G10").Value = 0.06904912275817471
$ws.Range("H10").Value = 0.05874797388748215
$ws.Range("I10").Value = 0.05182706490755814
$ws.Range("J10").Value = 0.0477462870105682
$ws.Range("K10").Value = 0.04617220239008156
$ws.Range("L10").Value = 0.04397883832909889
$ws.Range("M10").Value = 0.04181601043002973
$ws.Range("N10").Value = 0.04126264179257171
$ws.Range("O10").Value = 0.03974088021910053
$ws.Range("P10").Value = 0.03835819588876534
$ws.Range("Q10").Value = 0.0377406164495593
$ws.Range("R10").Value = 0.0368001503552299
$ws.Range("S10").Value = 0.0368001503552299
$ws.Range("T10").Value = 0.03642703282131078
$ws.Range("U10").Value = 0.03626125020393249
$ws.Range("V10").Value = 0.03583956568573277
$ws.Range("W10").Value = 0.03558098258102033
$ws.Range("X10").Value = 0.03548247742237016
$ws.Range("Y10").Value = 0.03524026357476553
$ws.Range("C11").Value = 0.5312497615814209
$ws.Range("E11").Value = 1727.073000699475
$ws.Range("F11").Value = 0.08601700353881347
$ws.Range("G11").Value = 0.06786659585851582
$ws.Range("H11").Value = 0.05769527257236647
$ws.Range("I11").Value = 0.054487598805918
$ws.Range("J11").Value = 0.04966698138800563
$ws.Range("K11").Value = 0.04528757105149364
$ws.Range("L11").Value = 0.04354265171496675
$ws.Range("M11").Value = 0.0410408828391178
$ws.Range("N11").Value = 0.03953801086422708
$ws.Range("O11").Value = 0.03852363971535211
$ws.Range("P11").Value = 0.03798951331632679
$ws.Range("Q11").Value = 0.03688574028108508
$ws.Range("R11").Value = 0.03577517091126691
$ws.Range("S11").Value = 0.03527650282062461
$ws.Range("T11").Value = 0.03500066483031491
$ws.Range("U11").Value = 0.03459427012250844
$ws.Range("V11").Value = 0.0342370161292235
$ws.Range("W11").Value = 0.03408210435227466
$ws.Range("X11").Value = 0.03398520827667384
$ws.Range("Y11").Value = 0.03366614036451217

Write-Host "Updated cells with new simulation run values"
